# Apply region_codes.xlsx edit: replace region rows 2-11 (North-West regions)
# with new rows 2-10 (Ural regions), shrinking the table by one row,
# adjust column widths, selection, and disable concurrent calc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (region_code, region_name) replacing the old 10 rows with 9 rows.
$codes = @(66, 45, 56, 59, 72, 86, 89, 74, 2)
$names = @(
    "Свердловская область",
    "Курганская область",
    "Оренбурская область",
    "Пермский край",
    "Тюменская область",
    "Ханты-мансийский АО - Югра",
    "Ямало-ненецкий АО",
    "Челябинская область",
    "Республика Башкортостан"
)

# Clear out the old used range below the header row (rows 2:11) first, since
# the new table is shorter by one row.
$ws.Range("A2:B11").Clear()

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $codes[$i]
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

# Column widths (widened to fit the longer Ural region names).
$ws.Columns.Item(1).ColumnWidth = 15.25
$ws.Columns.Item(2).ColumnWidth = 44.25

# Selection moves to B13 (below the shrunk table).
$ws.Range("B13").Select()

# Disable concurrent calculation.
$wb.ConcurrentCalc = $false
